# Updated symbol list on Sun Dec 11 22:46:18 UTC 2022 with GitHub Actions
#
# Applies the latest price-refresh pass to the "cryptos" sheet: most rows
# just get a refreshed Price (column D) value. A couple of coins
# (KickToken / BKEXToken) also swapped ranking positions, which moves
# their Coin name / Link / Volume(1h) label cells around too.
#
# Numeric-looking values must stay TEXT (the sheet stores Price as
# inline strings, not numbers), so they are entered with a leading
# apostrophe - exactly like a user forcing text entry in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $ws.Range($cellRef).Value = "'" + $text
}

# --- Column D (Price) refreshes -------------------------------------------
Set-TextValue "D2"  "286.05"
Set-TextValue "D3"  "21.19"
Set-TextValue "D4"  "6.454"
Set-TextValue "D5"  "0.06370"
Set-TextValue "D6"  "3.605"
Set-TextValue "D7"  "1.536"
Set-TextValue "D8"  "6.559"
Set-TextValue "D9"  "0.8216"
Set-TextValue "D10" "0.01403"
Set-TextValue "D11" "0.1679"
Set-TextValue "D12" "0.08691"
Set-TextValue "D13" "0.03670"
Set-TextValue "D14" "0.03215"
Set-TextValue "D15" "0.09190"
Set-TextValue "D16" "3.706"
Set-TextValue "D18" "0.04763"
Set-TextValue "D19" "0.006193"
Set-TextValue "D20" "0.006271"
Set-TextValue "D21" "0.001071"
Set-TextValue "D22" "0.0001601"
Set-TextValue "D23" "3.782"
Set-TextValue "D24" "2.270"
Set-TextValue "D25" "0.3358"
Set-TextValue "D26" "0.1263"
Set-TextValue "D40" "0.04775"

# --- Rows 41-43: KickToken / BKEXToken swap positions ----------------------
# Row 41 becomes BKEXToken (was KickToken)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1115"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 stays CEJI, but price/label refresh (no longer the 24h worst performer)
Set-TextValue "D42" "0.003451"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes KickToken (was BKEXToken), now flagged as 24h worst performer
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003582"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Remaining column D (Price) refreshes ----------------------------------
Set-TextValue "D45" "0.00006936"
Set-TextValue "D46" "0.00000000752"
Set-TextValue "D47" "1.003"
Set-TextValue "D48" "0.005459"
Set-TextValue "D49" "0.00001504"
Set-TextValue "D50" "0.01244"
